$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new student record as row 53
$ws.Range("A53").Value = 50
$ws.Range("B53").Value = "Chokkakula Kusuma"
$ws.Range("C53").Value = "24NAG0555_PC50"
$ws.Range("D53").Value = "Pass@word1"

# Add the mailto hyperlink for the password cell, matching the existing rows
$ws.Hyperlinks.Add($ws.Range("D53"), "mailto:Pass@word1") | Out-Null

# Re-apply the built-in Hyperlink cell style (Hyperlinks.Add nudges font formatting)
$ws.Range("D53").Style = "Hyperlink"

# Update the view so the new row is selected, like in the saved workbook
$ws.Range("A53").Select() | Out-Null

Write-Host "Added row 53 for Chokkakula Kusuma"
